$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 412.5
$ws.Range("I33").Value = 383.33334
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 383.33334
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -154.33334
$ws.Range("N33").Value = -958

$ws.Range("H74").Value = 3186771.8
$ws.Range("I74").Value = 4247787
$ws.Range("K74").Value = 4247787
$ws.Range("M74").Value = -4246851

$ws.Range("H77").Value = 3186771.8
$ws.Range("I77").Value = 4247787
$ws.Range("K77").Value = 21238935
$ws.Range("M77").Value = -21234255

$ws.Range("H121").Value = 1061
$ws.Range("J121").Value = 1235
$ws.Range("L121").Value = 3705
$ws.Range("N121").Value = -7199

$ws.Range("H131").Value = 3103
$ws.Range("I131").Value = 1600
$ws.Range("K131").Value = 4800
$ws.Range("M131").Value = 240

$ws.Range("H132").Value = 3158.2144
$ws.Range("I132").Value = 3158.2144
$ws.Range("K132").Value = 9474.643199999999
$ws.Range("M132").Value = -6944.643199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 15023.667
$ws.Range("I28").Value = 6487.5713
$ws.Range("J28").Value = 44900
$ws.Range("K28").Value = 6487.5713
$ws.Range("L28").Value = 44900
$ws.Range("M28").Value = -6295.5713
$ws.Range("N28").Value = -45284

$ws.Range("H32").Value = 6299.93
$ws.Range("I32").Value = 5191.946
$ws.Range("J32").Value = 13132.5
$ws.Range("K32").Value = 5191.946
$ws.Range("L32").Value = 13132.5
$ws.Range("M32").Value = -4904.946
$ws.Range("N32").Value = -13706.5

$ws.Range("H61").Value = 1404.8
$ws.Range("I61").Value = 1006
$ws.Range("K61").Value = 1006
$ws.Range("M61").Value = -794

$ws.Range("H99").Value = 15023.667
$ws.Range("I99").Value = 6487.5713
$ws.Range("J99").Value = 44900
$ws.Range("K99").Value = 6487.5713
$ws.Range("L99").Value = 44900
$ws.Range("M99").Value = -3492.5713
$ws.Range("N99").Value = -50890

$ws.Range("H110").Value = 1082.1666
$ws.Range("I110").Value = 900
$ws.Range("J110").Value = 1173.25
$ws.Range("K110").Value = 900
$ws.Range("L110").Value = 1173.25
$ws.Range("M110").Value = 1145
$ws.Range("N110").Value = -5263.25

$ws.Range("H122").Value = 1478.6
$ws.Range("I122").Value = 934.4
$ws.Range("J122").Value = 2022.8
$ws.Range("K122").Value = 2803.2
$ws.Range("L122").Value = 6068.4
$ws.Range("M122").Value = -353.1999999999998
$ws.Range("N122").Value = -10968.4

$ws.Range("H136").Value = 1404.8
$ws.Range("I136").Value = 1006
$ws.Range("K136").Value = 3018
$ws.Range("M136").Value = -468

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3784.0356
$ws.Range("I105").Value = 3087.7273
$ws.Range("J105").Value = 6337.1665
$ws.Range("K105").Value = 3087.7273
$ws.Range("L105").Value = 6337.1665
$ws.Range("M105").Value = -1340.7273
$ws.Range("N105").Value = -9831.166499999999

$ws.Range("H107").Value = 7537.7856
$ws.Range("I107").Value = 2632.9412
$ws.Range("J107").Value = 15118
$ws.Range("K107").Value = 2632.9412
$ws.Range("L107").Value = 15118
$ws.Range("M107").Value = -712.9412000000002
$ws.Range("N107").Value = -18958

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3018
$ws.Range("I132").Value = 712
$ws.Range("J132").Value = 3479.2
$ws.Range("K132").Value = 2136
$ws.Range("L132").Value = 10437.6
$ws.Range("M132").Value = 394
$ws.Range("N132").Value = -15497.6

$ws.Range("H134").Value = 50001692
$ws.Range("I134").Value = 1679.5
$ws.Range("J134").Value = 250001740
$ws.Range("K134").Value = 5038.5
$ws.Range("L134").Value = 750005220
$ws.Range("M134").Value = -2503.5
$ws.Range("N134").Value = -750010290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 2100
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -852
$ws.Range("N92").Value = -3996

$ws.Range("H110").Value = 10811.875
$ws.Range("I110").Value = 3831.6667
$ws.Range("J110").Value = 15000
$ws.Range("K110").Value = 11495.0001
$ws.Range("L110").Value = 45000
$ws.Range("M110").Value = -7405.000100000001
$ws.Range("N110").Value = -53180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4166.591
$ws.Range("I80").Value = 4336.5386
$ws.Range("J80").Value = 3921.111
$ws.Range("K80").Value = 4336.5386
$ws.Range("L80").Value = 3921.111
$ws.Range("M80").Value = -3338.5386
$ws.Range("N80").Value = -5917.111

$ws.Range("H83").Value = 4166.591
$ws.Range("I83").Value = 4336.5386
$ws.Range("J83").Value = 3921.111
$ws.Range("K83").Value = 21682.693
$ws.Range("L83").Value = 19605.555
$ws.Range("M83").Value = -16690.693
$ws.Range("N83").Value = -29589.555

$ws.Range("H136").Value = 7222.45
$ws.Range("J136").Value = 7222.45
$ws.Range("L136").Value = 21667.35
$ws.Range("N136").Value = -26767.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1090.129
$ws.Range("I16").Value = 910.5789
$ws.Range("K16").Value = 910.5789
$ws.Range("M16").Value = -740.5789

$ws.Range("H82").Value = 3116.818
$ws.Range("I82").Value = 2464.6667
$ws.Range("J82").Value = 6051.5
$ws.Range("K82").Value = 2464.6667
$ws.Range("L82").Value = 6051.5
$ws.Range("M82").Value = -2103.6667
$ws.Range("N82").Value = -6773.5

$ws.Range("H85").Value = 3116.818
$ws.Range("I85").Value = 2464.6667
$ws.Range("J85").Value = 6051.5
$ws.Range("K85").Value = 2464.6667
$ws.Range("L85").Value = 6051.5
$ws.Range("M85").Value = -1216.6667
$ws.Range("N85").Value = -8547.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4562.75
$ws.Range("J62").Value = 3925
$ws.Range("L62").Value = 3925
$ws.Range("N62").Value = -5173

$ws.Range("H65").Value = 4562.75
$ws.Range("J65").Value = 3925
$ws.Range("L65").Value = 19625
$ws.Range("N65").Value = -25865

$ws.Range("H81").Value = 2493.5
$ws.Range("I81").Value = 2192.2
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 4384.4
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -3323.4
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 2493.5
$ws.Range("I84").Value = 2192.2
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 21922
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -16618
$ws.Range("N84").Value = -50608
